# Commit changes to raw data and description of the data
#
# 1. On the "Data" sheet, the D and E columns are swapped in meaning:
#    D used to hold "eye-color" (text) and E used to hold "waist-size"
#    (number). They are re-labelled/re-populated so D holds "Waist-size"
#    (number) and E holds "Eye-color" (text).
# 2. On the "Codebook" sheet, two new rows document the "waist-size" and
#    "eye-color" variables.
# 3. Minor view-state touch-ups (selection, page setup) follow what Excel
#    would naturally record after making the edits above.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Data" sheet: re-label headers and swap/repopulate columns D and E
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

$data.Cells.Item(1, 4).Value() = "Waist-size"
$data.Cells.Item(1, 5).Value() = "Eye-color"

$waist = @{
    2  = 32
    3  = 45
    4  = 40
    5  = 26
    6  = 27
    7  = 55
    8  = 42
    9  = 25
    10 = 34
    11 = 55
    12 = 43
    13 = 34
    14 = 30
    15 = 24
}
$eye = @{
    2  = "blue"
    3  = "black"
    4  = "green"
    5  = "brown"
    6  = "blue"
    7  = "green"
    8  = "black"
    9  = "black"
    10 = "brown"
    11 = "green"
    12 = "blue"
    13 = "brown"
    14 = "green"
    15 = "brown"
}

foreach ($r in 2..15) {
    $data.Cells.Item($r, 4).Value() = $waist[$r]
    $data.Cells.Item($r, 5).Value() = $eye[$r]
}

# Give the sheet the print-portrait page setup Excel stamps on after the
# edit session.
$data.PageSetup.Orientation() = 1

# ---------------------------------------------------------------------
# 2. "Codebook" sheet: document the two new/changed variables
# ---------------------------------------------------------------------
$codebook = $wb.Worksheets.Item("Codebook")

$codebook.Cells.Item(5, 1).Value() = "waist-size"
$codebook.Cells.Item(5, 2).Value() = "size of the waist in inches"
$codebook.Cells.Item(5, 3).Value() = "numeric value >0 or NA"

$codebook.Cells.Item(6, 1).Value() = "eye-color"
$codebook.Cells.Item(6, 2).Value() = "color of the individuals eye"
$codebook.Cells.Item(6, 3).Value() = "character variable"

# ---------------------------------------------------------------------
# 3. View state: leave selection on Codebook!C6, then Data!F2 so "Data"
#    ends up the active (tab-selected) sheet, matching the saved file.
# ---------------------------------------------------------------------
$codebook.Range("C6").Select()
$data.Range("F2").Select()
